# Update New Orleans xlsx: add a "State" column to hotel_info (right after
# Hotel_Name, before City) and reorder the sheets so review_info comes
# before hotel_info.

$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new "State" column after "Hotel_Name" (column B) in hotel_info.
$hotelSheet.Columns("C").Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Move review_info so it becomes the first sheet (before hotel_info).
$reviewSheet.Move($hotelSheet)
